$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Default column width hint (baseColWidth="10" in the saved XML).
$ws.StandardWidth = 10

# Insert a new first column to hold the "Version" field; everything else
# (Code / Description / Definition) shifts one column to the right.
$ws.Columns.Item(1).Insert()

# Header
$ws.Range("A1").Value = "Version"

# Data rows: every existing data row gets "1.0" as its version. Writing the
# literal text "1.0" via .Value would be auto-coerced to the number 1, so
# instead write it as a formula that evaluates to the text "1.0" and then
# flatten it down to a plain value in place (copy / paste-special-values) -
# this keeps the cell a genuine text cell without touching cell formatting.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Formula = '="1.0"'
}
$versionRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item($lastRow, 1))
$versionRange.Copy()
$versionRange.PasteSpecial(-4163)
$excel.CutCopyMode = $false
